# Weekly update: insert this week's new price records at the top of the
# data block (row 577) and push the previously-top rows down by 4, which
# is how this particular logging sheet keeps its most-recent-week-first
# ordering. Excel's native row insert takes care of shifting the existing
# 577:607 block down to 581:611 (formatting, including the date style on
# column D, comes along for free), so we only need to populate the four
# freshly inserted rows with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 577 down by inserting 4 fresh rows above it.
$ws.Rows("577:580").Insert()

# Fixed/common columns shared by every record in this block.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"
$unidad      = "$/bins (400 kilos)"
$origen      = "Provincia de Limarí"
$kgUnidad    = 400

function Set-Record {
    param($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-Record 577 44516 "Navel Late" "Primera" 24 170000 180000 175000 438
Set-Record 578 44516 "Navel Late" "Segunda" 20 150000 160000 155000 388
Set-Record 579 44516 "Valencia"   "Primera" 20 160000 170000 165000 412
Set-Record 580 44516 "Valencia"   "Segunda" 16 130000 140000 135000 338
